# Daily attendance processing - reorder "Recorded By" (column G) contributor
# lists. Each multi-contributor cell has its comma-separated entries
# reversed in place, except for the literal "admin@admin.com, System"
# combination which is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    if ($val -notmatch ",") {
        continue
    }

    if ($val -eq "admin@admin.com, System") {
        continue
    }

    $parts = $val -split ", "
    $reversed = $parts[($parts.Count - 1)..0]
    $newVal = $reversed -join ", "

    $cell.Value = $newVal
}
